# AltriaUS update to code
# Adds a new "Spacing" worksheet (with category Low/High/Size bucket data)
# after the existing "KPI" sheet, makes it the active sheet, and tweaks
# a couple of view/selection details on the KPI sheet.

$wb = $excel.ActiveWorkbook
$kpi = $wb.Worksheets.Item("KPI")

# --- Tweak the KPI sheet's view/selection ------------------------------
[void]$kpi.Range("G3").Select()

# --- Create the new "Spacing" sheet right after "KPI" ------------------
$ws = $wb.Worksheets.Add($null, $kpi)
$ws.Name = "Spacing"

# Header row
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Value = "Low"
$ws.Range("C1").Value = "High"
$ws.Range("D1").Value = "Size"

# Category spacing/bucket data
$data = @(
  @("Cigarettes", 1, 5, 1),
  @("Cigarettes", 6, 9, 2),
  @("Cigarettes", 10, 14, 3),
  @("Cigarettes", 15, 19, 4),
  @("Smokeless", 1, 4, 1),
  @("Smokeless", 5, 8, 2),
  @("Smokeless", 9, 12, 3),
  @("Smokeless", 13, 16, 4),
  @("Cigars", 1, 5, 1),
  @("Cigars", 6, 9, 2),
  @("Cigars", 10, 14, 3),
  @("Cigars", 15, 19, 4),
  @("Vapor", 1, 5, 1),
  @("Vapor", 6, 9, 2),
  @("Vapor", 10, 14, 3),
  @("Vapor", 15, 19, 4)
)

$row = 2
foreach ($d in $data) {
    $ws.Cells.Item($row, 1).Value = $d[0]
    $ws.Cells.Item($row, 2).Value = $d[1]
    $ws.Cells.Item($row, 3).Value = $d[2]
    $ws.Cells.Item($row, 4).Value = $d[3]
    $row++
}

# Column width for the new sheet (single column group, narrow width)
[void]$ws.Columns("A:D").EntireColumn

# --- Column width adjustments on the KPI sheet --------------------------
$kpi.Range("A1").EntireColumn.ColumnWidth = 29.8095238095238
$kpi.Range("B1").EntireColumn.ColumnWidth = 21.0340136054422
$kpi.Range("C1:AMK1").EntireColumn.ColumnWidth = 7.67142857142857

# --- Make "Spacing" the active sheet/tab (last, so it ends up selected) -
$ws.Activate()

Write-Host "Spacing sheet added after KPI; workbook updated."
